$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2608.2307
$ws.Range("J17").Value = 2608.2307
$ws.Range("L17").Value = 7824.6921
$ws.Range("N17").Value = -8160.6921
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()
$ws.Range("H80").Value = 912.5
$ws.Range("I80").Value = 650
$ws.Range("J80").Value = 1000
$ws.Range("K80").Value = 1950
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -952
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 912.5
$ws.Range("I83").Value = 650
$ws.Range("J83").Value = 1000
$ws.Range("K83").Value = 5850
$ws.Range("L83").Value = 9000
$ws.Range("M83").Value = -858
$ws.Range("N83").Value = -18984
$ws.Range("H127").Value = 1237.5
$ws.Range("I127").Value = 1237.5
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 3712.5
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 1247.5
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1008
$ws.Range("I2").Value = 1008
$ws.Range("K2").Value = 1008
$ws.Range("M2").Value = -895
$ws.Range("H14").Value = 3612.25
$ws.Range("I14").Value = 3225
$ws.Range("J14").Value = 3999.5
$ws.Range("K14").Value = 3225
$ws.Range("L14").Value = 3999.5
$ws.Range("M14").Value = -3050
$ws.Range("N14").Value = -4349.5
$ws.Range("H43").Value = 12499984
$ws.Range("J43").Value = 9999968
$ws.Range("L43").Value = 9999968
$ws.Range("N43").Value = -10000594
$ws.Range("H97").Value = 845.7059
$ws.Range("I97").Value = 858.3077
$ws.Range("K97").Value = 858.3077
$ws.Range("M97").Value = -362.3077
$ws.Range("H116").Value = 1008
$ws.Range("I116").Value = 1008
$ws.Range("K116").Value = 1008
$ws.Range("M116").Value = 1286
$ws.Range("H132").Value = 8954.6
$ws.Range("I132").Value = 8954.6
$ws.Range("K132").Value = 26863.8
$ws.Range("M132").Value = -24333.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1008
$ws.Range("I3").Value = 1008
$ws.Range("K3").Value = 1008
$ws.Range("M3").Value = -894

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 12913.875
$ws.Range("I7").Value = 16920.166
$ws.Range("K7").Value = 16920.166
$ws.Range("M7").Value = -16807.166
$ws.Range("H31").Value = 8945
$ws.Range("I31").Value = 3142.5
$ws.Range("K31").Value = 3142.5
$ws.Range("M31").Value = -2847.5
$ws.Range("H34").Value = 8945
$ws.Range("I34").Value = 3142.5
$ws.Range("K34").Value = 3142.5
$ws.Range("M34").Value = -2940.5
$ws.Range("H41").Value = 34062.793
$ws.Range("J41").Value = 34456.824
$ws.Range("L41").Value = 34456.824
$ws.Range("N41").Value = -35312.824
$ws.Range("H47").Value = 22277.5
$ws.Range("I47").Value = 22277.5
$ws.Range("K47").Value = 22277.5
$ws.Range("M47").Value = -21711.5
$ws.Range("H59").Value = 329280000
$ws.Range("I59").Value = 86900
$ws.Range("K59").Value = 86900
$ws.Range("M59").Value = -85755
$ws.Range("H62").Value = 8055.778
$ws.Range("I62").Value = 7099.8
$ws.Range("J62").Value = 9250.75
$ws.Range("K62").Value = 7099.8
$ws.Range("L62").Value = 9250.75
$ws.Range("M62").Value = -6475.8
$ws.Range("N62").Value = -10498.75
$ws.Range("H65").Value = 8055.778
$ws.Range("I65").Value = 7099.8
$ws.Range("J65").Value = 9250.75
$ws.Range("K65").Value = 35499
$ws.Range("L65").Value = 46253.75
$ws.Range("M65").Value = -32379
$ws.Range("N65").Value = -52493.75
$ws.Range("H76").Value = 6250
$ws.Range("I76").Value = 6250
$ws.Range("K76").Value = 6250
$ws.Range("M76").Value = -5935
$ws.Range("H79").Value = 6250
$ws.Range("I79").Value = 6250
$ws.Range("K79").Value = 6250
$ws.Range("M79").Value = -5158
$ws.Range("H104").Value = 40284.5
$ws.Range("J104").Value = 40284.5
$ws.Range("L104").Value = 40284.5
$ws.Range("N104").Value = -45526.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 14699174
$ws.Range("I11").Value = 10593000
$ws.Range("K11").Value = 10593000
$ws.Range("M11").Value = -10592861
$ws.Range("H15").Value = 20008.5
$ws.Range("J15").Value = 20008.5
$ws.Range("L15").Value = 20008.5
$ws.Range("N15").Value = -20584.5
$ws.Range("H81").Value = 20008.5
$ws.Range("J81").Value = 20008.5
$ws.Range("L81").Value = 20008.5
$ws.Range("N81").Value = -22004.5
$ws.Range("H84").Value = 20008.5
$ws.Range("J84").Value = 20008.5
$ws.Range("L84").Value = 60025.5
$ws.Range("N84").Value = -70009.5
$ws.Range("H97").Value = 425.85
$ws.Range("I97").Value = 346.92307
$ws.Range("J97").Value = 572.4286
$ws.Range("K97").Value = 346.92307
$ws.Range("L97").Value = 572.4286
$ws.Range("M97").Value = 149.07693
$ws.Range("N97").Value = -1564.4286
$ws.Range("H126").Value = 6132.4
$ws.Range("I126").Value = 6912
$ws.Range("K126").Value = 20736
$ws.Range("M126").Value = -18266

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 11750
$ws.Range("I25").Value = 11750
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 11750
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -11520
$ws.Range("N25").ClearContents()
$ws.Range("H46").Value = 6933.3335
$ws.Range("J46").Value = 8000
$ws.Range("L46").Value = 8000
$ws.Range("N46").Value = -8376
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6990.1
$ws.Range("I62").Value = 2100.5
$ws.Range("K62").Value = 2100.5
$ws.Range("M62").Value = -1476.5
$ws.Range("H65").Value = 6990.1
$ws.Range("I65").Value = 2100.5
$ws.Range("K65").Value = 10502.5
$ws.Range("M65").Value = -7382.5
